$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 311 entirely ("「時は剣のごとし。あなたが割かなければあなたが割かれる」"),
# shifting all subsequent rows up by one.
$ws.Rows.Item(311).Delete()
